$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 383800.47
$ws.Range("J17").Value = 383800.47
$ws.Range("L17").Value = 1151401.41
$ws.Range("N17").Value = -1151737.41
# Row 40
$ws.Range("H40").Value = 5552.2
$ws.Range("I40").Value = 4602.933
$ws.Range("K40").Value = 4602.933
$ws.Range("M40").Value = -4427.933
# Row 74
$ws.Range("H74").Value = 4392.2144
$ws.Range("I74").Value = 3436.375
$ws.Range("K74").Value = 3436.375
$ws.Range("M74").Value = -2500.375
# Row 75
$ws.Range("H75").Value = 37657
$ws.Range("J75").Value = 37657
$ws.Range("L75").Value = 37657
$ws.Range("N75").Value = -39529
# Row 77
$ws.Range("H77").Value = 4392.2144
$ws.Range("I77").Value = 3436.375
$ws.Range("K77").Value = 17181.875
$ws.Range("M77").Value = -12501.875
# Row 78
$ws.Range("H78").Value = 37657
$ws.Range("J78").Value = 37657
$ws.Range("L78").Value = 112971
$ws.Range("N78").Value = -122331
# Row 80
$ws.Range("H80").Value = 1619.12
$ws.Range("I80").Value = 825.3333
$ws.Range("J80").Value = 2809.8
$ws.Range("K80").Value = 2475.9999
$ws.Range("L80").Value = 8429.400000000001
$ws.Range("M80").Value = -1477.9999
$ws.Range("N80").Value = -10425.4
# Row 83
$ws.Range("H83").Value = 1619.12
$ws.Range("I83").Value = 825.3333
$ws.Range("J83").Value = 2809.8
$ws.Range("K83").Value = 7427.9997
$ws.Range("L83").Value = 25288.2
$ws.Range("M83").Value = -2435.9997
$ws.Range("N83").Value = -35272.2
# Row 116
$ws.Range("H116").Value = 4505.222
$ws.Range("I116").Value = 3899.5
$ws.Range("J116").Value = 4989.8
$ws.Range("K116").Value = 3899.5
$ws.Range("L116").Value = 4989.8
$ws.Range("M116").Value = -457.5
$ws.Range("N116").Value = -11873.8
# Row 132
$ws.Range("H132").Value = 3463.1936
$ws.Range("I132").Value = 1707.0714
$ws.Range("K132").Value = 5121.2142
$ws.Range("M132").Value = -2591.2142
# Row 137
$ws.Range("H137").Value = 1242.9
$ws.Range("I137").Value = 1718
$ws.Range("J137").Value = 1124.125
$ws.Range("K137").Value = 5154
$ws.Range("L137").Value = 3372.375
$ws.Range("M137").Value = -2604
$ws.Range("N137").Value = -8472.375

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 4247.569
$ws.Range("I32").Value = 3764.5557
$ws.Range("K32").Value = 3764.5557
$ws.Range("M32").Value = -3477.5557
# Row 45
$ws.Range("H45").Value = 6943.12
$ws.Range("I45").Value = 9615.385
$ws.Range("K45").Value = 9615.385
$ws.Range("M45").Value = -9238.385
# Row 61
$ws.Range("H61").Value = 2465.7693
$ws.Range("I61").Value = 1342.7778
$ws.Range("K61").Value = 1342.7778
$ws.Range("M61").Value = -1130.7778
# Row 136
$ws.Range("H136").Value = 2465.7693
$ws.Range("I136").Value = 1342.7778
$ws.Range("K136").Value = 4028.3334
$ws.Range("M136").Value = -1478.3334

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 4623.7446
$ws.Range("I20").Value = 4754.161
$ws.Range("J20").Value = 4371.0625
$ws.Range("K20").Value = 4754.161
$ws.Range("L20").Value = 4371.0625
$ws.Range("M20").Value = -4507.161
$ws.Range("N20").Value = -4865.0625
# Row 86
$ws.Range("H86").Value = 2401
$ws.Range("I86").Value = 2449.2964
$ws.Range("J86").Value = 1966.3334
$ws.Range("K86").Value = 2449.2964
$ws.Range("L86").Value = 1966.3334
$ws.Range("M86").Value = -1326.2964
$ws.Range("N86").Value = -4212.3334
# Row 89
$ws.Range("H89").Value = 2401
$ws.Range("I89").Value = 2449.2964
$ws.Range("J89").Value = 1966.3334
$ws.Range("K89").Value = 12246.482
$ws.Range("L89").Value = 9831.666999999999
$ws.Range("M89").Value = -6630.482
$ws.Range("N89").Value = -21063.667
# Row 94
$ws.Range("H94").Value = 2039.8
$ws.Range("I94").Value = 2039.8
$ws.Range("K94").Value = 2039.8
$ws.Range("M94").Value = -1588.8
# Row 134
$ws.Range("H134").Value = 1178
$ws.Range("I134").Value = 1178
$ws.Range("K134").Value = 3534
$ws.Range("M134").Value = -999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 4805.7334
$ws.Range("I31").Value = 2661.889
$ws.Range("J31").Value = 8021.5
$ws.Range("K31").Value = 2661.889
$ws.Range("L31").Value = 8021.5
$ws.Range("M31").Value = -2366.889
$ws.Range("N31").Value = -8611.5
# Row 34
$ws.Range("H34").Value = 4805.7334
$ws.Range("I34").Value = 2661.889
$ws.Range("J34").Value = 8021.5
$ws.Range("K34").Value = 2661.889
$ws.Range("L34").Value = 8021.5
$ws.Range("M34").Value = -2459.889
$ws.Range("N34").Value = -8425.5
# Row 93
$ws.Range("H93").Value = 15772.429
$ws.Range("I93").Value = 9902.200000000001
$ws.Range("K93").Value = 9902.200000000001
$ws.Range("M93").Value = -8030.200000000001
# Row 132
$ws.Range("H132").Value = 5527
$ws.Range("I132").Value = 5527
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 16581
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -14051
$ws.Range("N132").ClearContents()
# Row 134
$ws.Range("H134").Value = 2205.9644
$ws.Range("I134").Value = 1270.68
$ws.Range("J134").Value = 10000
$ws.Range("K134").Value = 3812.04
$ws.Range("L134").Value = 30000
$ws.Range("M134").Value = -1277.04
$ws.Range("N134").Value = -35070

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 332.76923
$ws.Range("I2").Value = 519.1539
$ws.Range("J2").Value = 146.38461
$ws.Range("K2").Value = 3114.9234
$ws.Range("L2").Value = 878.3076600000001
$ws.Range("M2").Value = -3001.9234
$ws.Range("N2").Value = -1104.30766
# Row 38
$ws.Range("H38").Value = 68.61539
$ws.Range("I38").Value = 31.375
$ws.Range("J38").Value = 128.2
$ws.Range("K38").Value = 94.125
$ws.Range("L38").Value = 384.6
$ws.Range("M38").Value = 252.875
$ws.Range("N38").Value = -1078.6
# Row 81
$ws.Range("H81").Value = 55564830
$ws.Range("J81").Value = 100010500
$ws.Range("L81").Value = 300031500
$ws.Range("N81").Value = -300033746
# Row 84
$ws.Range("H84").Value = 55564830
$ws.Range("J84").Value = 100010500
$ws.Range("L84").Value = 900094500
$ws.Range("N84").Value = -900105732
# Row 86
$ws.Range("H86").Value = 786.13635
$ws.Range("I86").Value = 651.875
$ws.Range("J86").Value = 862.8570999999999
$ws.Range("K86").Value = 1955.625
$ws.Range("L86").Value = 2588.5713
$ws.Range("M86").Value = -769.625
$ws.Range("N86").Value = -4960.5713
# Row 89
$ws.Range("H89").Value = 786.13635
$ws.Range("I89").Value = 651.875
$ws.Range("J89").Value = 862.8570999999999
$ws.Range("K89").Value = 5866.875
$ws.Range("L89").Value = 7765.7139
$ws.Range("M89").Value = 61.125
$ws.Range("N89").Value = -19621.7139
# Row 103
$ws.Range("H103").Value = 657.44446
$ws.Range("I103").Value = 497.66666
$ws.Range("J103").Value = 737.3333
$ws.Range("K103").Value = 1492.99998
$ws.Range("L103").Value = 2211.9999
$ws.Range("M103").Value = -613.9999800000001
$ws.Range("N103").Value = -3969.9999
# Row 120
$ws.Range("H120").Value = 18199.8
$ws.Range("I120").Value = 17749.75
$ws.Range("J120").Value = 20000
$ws.Range("K120").Value = 53249.25
$ws.Range("L120").Value = 60000
$ws.Range("M120").Value = -48411.25
$ws.Range("N120").Value = -69676
# Row 129
$ws.Range("H129").Value = 93449.5
$ws.Range("I129").Value = 500524.75
$ws.Range("J129").Value = 2988.3333
$ws.Range("K129").Value = 1501574.25
$ws.Range("L129").Value = 8964.999899999999
$ws.Range("M129").Value = -1496574.25
$ws.Range("N129").Value = -18964.9999
# Row 131
$ws.Range("H131").Value = 5631.3794
$ws.Range("J131").Value = 2376
$ws.Range("L131").Value = 7128
$ws.Range("N131").Value = -17208

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 7
$ws.Range("H7").Value = 5258693.5
$ws.Range("I7").Value = 700
$ws.Range("J7").Value = 27167000
$ws.Range("K7").Value = 700
$ws.Range("L7").Value = 27167000
$ws.Range("M7").Value = -588
$ws.Range("N7").Value = -27167224
# Row 8
$ws.Range("H8").Value = 5258693.5
$ws.Range("I8").Value = 700
$ws.Range("J8").Value = 27167000
$ws.Range("K8").Value = 700
$ws.Range("L8").Value = 27167000
$ws.Range("M8").Value = -561
$ws.Range("N8").Value = -27167278
# Row 29
$ws.Range("H29").Value = 6900
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()
# Row 80
$ws.Range("H80").Value = 3683.9
$ws.Range("I80").Value = 2981.125
$ws.Range("K80").Value = 2981.125
$ws.Range("M80").Value = -1983.125
# Row 83
$ws.Range("H83").Value = 3683.9
$ws.Range("I83").Value = 2981.125
$ws.Range("K83").Value = 14905.625
$ws.Range("M83").Value = -9913.625
# Row 97
$ws.Range("H97").Value = 66236.25
$ws.Range("I97").Value = 86874.44500000001
$ws.Range("K97").Value = 86874.44500000001
$ws.Range("M97").Value = -86378.44500000001
# Row 102
$ws.Range("H102").Value = 1165.1428
$ws.Range("I102").Value = 1212.1875
$ws.Range("K102").Value = 1212.1875
$ws.Range("M102").Value = 409.8125
# Row 113
$ws.Range("H113").Value = 2915.3572
$ws.Range("I113").Value = 2113.889
$ws.Range("K113").Value = 2113.889
$ws.Range("M113").Value = 56.11099999999988
# Row 122
$ws.Range("H122").Value = 5081
$ws.Range("I122").Value = 4097.4
$ws.Range("J122").Value = 9999
$ws.Range("K122").Value = 12292.2
$ws.Range("L122").Value = 29997
$ws.Range("M122").Value = -9842.199999999999
$ws.Range("N122").Value = -34897

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 551.4167
$ws.Range("I16").Value = 551.4167
$ws.Range("K16").Value = 551.4167
$ws.Range("M16").Value = -381.4167
# Row 55
$ws.Range("H55").Value = 2773.6072
$ws.Range("I55").Value = 2786.0667
$ws.Range("J55").Value = 2759.2307
$ws.Range("K55").Value = 2786.0667
$ws.Range("L55").Value = 2759.2307
$ws.Range("M55").Value = -2613.0667
$ws.Range("N55").Value = -3105.2307
# Row 68
$ws.Range("H68").Value = 1999
$ws.Range("I68").Value = 1999
$ws.Range("K68").Value = 1999
$ws.Range("M68").Value = -1250
# Row 71
$ws.Range("H71").Value = 1999
$ws.Range("I71").Value = 1999
$ws.Range("K71").Value = 9995
$ws.Range("M71").Value = -6251
# Row 82
$ws.Range("H82").Value = 3106.8333
$ws.Range("I82").Value = 2880.3333
$ws.Range("J82").Value = 3333.3333
$ws.Range("K82").Value = 2880.3333
$ws.Range("L82").Value = 3333.3333
$ws.Range("M82").Value = -2519.3333
$ws.Range("N82").Value = -4055.3333
# Row 85
$ws.Range("H85").Value = 3106.8333
$ws.Range("I85").Value = 2880.3333
$ws.Range("J85").Value = 3333.3333
$ws.Range("K85").Value = 2880.3333
$ws.Range("L85").Value = 3333.3333
$ws.Range("M85").Value = -1632.3333
$ws.Range("N85").Value = -5829.3333
# Row 93
$ws.Range("H93").Value = 15248.346
$ws.Range("I93").Value = 2184.652
$ws.Range("J93").Value = 115403.336
$ws.Range("K93").Value = 2184.652
$ws.Range("L93").Value = 115403.336
$ws.Range("M93").Value = -936.652
$ws.Range("N93").Value = -117899.336
# Row 100
$ws.Range("H100").Value = 29944.285
$ws.Range("I100").Value = 6538.4
$ws.Range("J100").Value = 42947.555
$ws.Range("K100").Value = 6538.4
$ws.Range("L100").Value = 42947.555
$ws.Range("M100").Value = -5997.4
$ws.Range("N100").Value = -44029.555

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 96
$ws.Range("H96").Value = 1747.4
$ws.Range("I96").Value = 1698.6
$ws.Range("J96").Value = 1796.2
$ws.Range("K96").Value = 1698.6
$ws.Range("L96").Value = 1796.2
$ws.Range("M96").Value = -325.5999999999999
$ws.Range("N96").Value = -4542.2
# Row 122
$ws.Range("H122").Value = 1757.7273
$ws.Range("I122").Value = 1757.7273
$ws.Range("K122").Value = 5273.1819
$ws.Range("M122").Value = -2823.1819
# Row 132
$ws.Range("H132").Value = 2615.5151
$ws.Range("I132").Value = 2461.8386
$ws.Range("J132").Value = 4997.5
$ws.Range("K132").Value = 7385.5158
$ws.Range("L132").Value = 14992.5
$ws.Range("M132").Value = -4855.5158
$ws.Range("N132").Value = -20052.5
# Row 136
$ws.Range("H136").Value = 2314.4285
$ws.Range("I136").Value = 2633.6667
$ws.Range("K136").Value = 7901.000100000001
$ws.Range("M136").Value = -5351.000100000001
